# Updated cryptos list on Thu Oct 12 23:30:02 UTC 2023 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures for the
# cryptos table on the active sheet, and re-syncs a handful of Coin/Link rows
# (B/C columns) whose ranking order changed, per the upstream data source.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.831.22'
$ws.Range('E2').Value = '  -0.02%  '
$ws.Range('D3').Value = '1.542.50'
$ws.Range('E3').Value = '  -1.61%  '
$ws.Range('E4').Value = '  +0.31%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '206.05'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.33%  '
$ws.Range('E6').Value = '  -0.86%  '
$ws.Range('E8').Value = '  -0.54%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '21.40'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -2.77%  '
$ws.Range('E10').Value = '  -0.88%  '
$ws.Range('E11').Value = '  -1.11%  '
$ws.Range('D12').Value = '1.761.89'
$ws.Range('E12').Value = '  -1.57%  '
$ws.Range('D13').Value = '1.548.15'
$ws.Range('E13').Value = '  -1.30%  '
$ws.Range('E14').Value = '  -1.58%  '
$ws.Range('E15').Value = '  -1.16%  '
$ws.Range('D16').Value = '26.828.75'
$ws.Range('E16').Value = '  -0.05%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '61.25'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -0.45%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '215.00'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +0.13%  '
$ws.Range('E19').Value = '  -2.58%  '
$ws.Range('E20').Value = '  +0.51%  '
$ws.Range('E21').Value = '  +0.24%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '9.15'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -1.61%  '
$ws.Range('E24').Value = '  -2.70%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '153.05'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.36%  '
$ws.Range('E26').Value = '  -1.89%  '
$ws.Range('E27').Value = '  -0.99%  '
$ws.Range('E28').Value = '  +0.25%  '
$ws.Range('E29').Value = '  -0.74%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.10'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -1.61%  '
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.0458'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -2.07%  '
$ws.Range('E32').Value = '  +1.35%  '
$ws.Range('D33').Value = '1.369.52'
$ws.Range('E33').Value = '  -2.41%  '
$ws.Range('E34').Value = '  +0.51%  '
$ws.Range('E35').Value = '  -1.31%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.960'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +2.53%  '
$ws.Range('E37').Value = '  +0.10%  '
$ws.Range('E38').Value = '  +1.36%  '
$ws.Range('E39').Value = '  -1.67%  '
$ws.Range('E40').Value = '  +8.96%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.806'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -1.10%  '
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.00'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +0.22%  '
$ws.Range('B43').Value = 'WEMIXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.991'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -0.04%  '
$ws.Range('B44').Value = 'MXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.22'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +1.53%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '63.19'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -0.27%  '
$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.73'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -3.99%  '
$ws.Range('B47').Value = 'RocketPoolETH'
$ws.Range('C47').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D47').Value = '1.676.47'
$ws.Range('E47').Value = '  -1.52%  '
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '84.16'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -2.34%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0511'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +3.79%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₇0970'
$ws.Range('E50').Value = '  -1.21%  '
$ws.Range('E51').Value = '  +0.20%  '
